$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Week 49 header (text, matches style of existing week-number headers)
$ws.Range("AZ1").Value = "'49"

# Numeric data cells for week 49 (column AZ), plus two sparse-row backfills
$ws.Range("AZ2").Value = 0
$ws.Range("AZ5").Value = 0
$ws.Range("AZ6").Value = 15
$ws.Range("AZ7").Value = 1
$ws.Range("AZ8").Value = 8
$ws.Range("AZ9").Value = 0
$ws.Range("AZ10").Value = 0
$ws.Range("AZ11").Value = 0
$ws.Range("AZ12").Value = 0
$ws.Range("AZ14").Value = 0
$ws.Range("AZ16").Value = 0
$ws.Range("AZ17").Value = 0
$ws.Range("AZ23").Value = 0
$ws.Range("AZ25").Value = 2
$ws.Range("AZ26").Value = 0
$ws.Range("AZ28").Value = 7
$ws.Range("AZ29").Value = 1
$ws.Range("AZ31").Value = 0
$ws.Range("AZ35").Value = 6
$ws.Range("AY36").Value = 0
$ws.Range("AZ36").Value = 0
$ws.Range("AZ37").Value = 0
$ws.Range("AZ38").Value = 0
$ws.Range("AZ41").Value = 0
$ws.Range("AZ42").Value = 0
$ws.Range("X43").Value = 0
$ws.Range("AA43").Value = 0
$ws.Range("AZ43").Value = 0
$ws.Range("AZ44").Value = 0
$ws.Range("AZ45").Value = 0
$ws.Range("AZ46").Value = 0
$ws.Range("AZ47").Value = 0
$ws.Range("AZ48").Value = 0
$ws.Range("AZ49").Value = 0
$ws.Range("AZ50").Value = 0
$ws.Range("AZ51").Value = 0
$ws.Range("AZ54").Value = 0
$ws.Range("AZ55").Value = 0
$ws.Range("AZ56").Value = 0
$ws.Range("AZ57").Value = 0
$ws.Range("AZ58").Value = 0
$ws.Range("AZ59").Value = 0
